{"js": "// Replace the two-digit multiplication expressions in the practice table\n// with the newly generated set of problems. Each old expression is unique\n// in the document, so a simple search + replace on each pair is safe.\nconst replacements = [\n  [\"21\u00d719=\", \"18\u00d791=\"],\n  [\"92\u00d711=\", \"78\u00d758=\"],\n  [\"53\u00d795=\", \"18\u00d790=\"],\n  [\"86\u00d725=\", \"42\u00d797=\"],\n  [\"42\u00d740=\", \"49\u00d777=\"],\n  [\"31\u00d727=\", \"62\u00d750=\"],\n  [\"40\u00d776=\", \"21\u00d797=\"],\n  [\"68\u00d711=\", \"18\u00d779=\"],\n  [\"88\u00d714=\", \"84\u00d723=\"],\n  [\"91\u00d719=\", \"86\u00d712=\"],\n  [\"29\u00d785=\", \"85\u00d734=\"],\n  [\"47\u00d724=\", \"24\u00d721=\"],\n  [\"22\u00d762=\", \"50\u00d756=\"],\n  [\"28\u00d723=\", \"38\u00d788=\"],\n  [\"62\u00d738=\", \"64\u00d756=\"],\n  [\"52\u00d790=\", \"76\u00d711=\"],\n  [\"47\u00d723=\", \"47\u00d754=\"],\n  [\"39\u00d758=\", \"68\u00d755=\"],\n  [\"33\u00d780=\", \"58\u00d781=\"],\n  [\"30\u00d768=\", \"37\u00d722=\"],\n  [\"36\u00d741=\", \"86\u00d729=\"],\n  [\"21\u00d720=\", \"13\u00d763=\"],\n  [\"76\u00d792=\", \"69\u00d725=\"],\n  [\"98\u00d732=\", \"77\u00d725=\"],\n  [\"88\u00d772=\", \"66\u00d755=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the two-digit multiplication expressions in the practice table\n# with the newly generated set of problems. Each old expression is unique\n# in the document, so Find/Replace on each pair (replace-one) is safe.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"21\u00d719=\", \"18\u00d791=\"),\n  @(\"92\u00d711=\", \"78\u00d758=\"),\n  @(\"53\u00d795=\", \"18\u00d790=\"),\n  @(\"86\u00d725=\", \"42\u00d797=\"),\n  @(\"42\u00d740=\", \"49\u00d777=\"),\n  @(\"31\u00d727=\", \"62\u00d750=\"),\n  @(\"40\u00d776=\", \"21\u00d797=\"),\n  @(\"68\u00d711=\", \"18\u00d779=\"),\n  @(\"88\u00d714=\", \"84\u00d723=\"),\n  @(\"91\u00d719=\", \"86\u00d712=\"),\n  @(\"29\u00d785=\", \"85\u00d734=\"),\n  @(\"47\u00d724=\", \"24\u00d721=\"),\n  @(\"22\u00d762=\", \"50\u00d756=\"),\n  @(\"28\u00d723=\", \"38\u00d788=\"),\n  @(\"62\u00d738=\", \"64\u00d756=\"),\n  @(\"52\u00d790=\", \"76\u00d711=\"),\n  @(\"47\u00d723=\", \"47\u00d754=\"),\n  @(\"39\u00d758=\", \"68\u00d755=\"),\n  @(\"33\u00d780=\", \"58\u00d781=\"),\n  @(\"30\u00d768=\", \"37\u00d722=\"),\n  @(\"36\u00d741=\", \"86\u00d729=\"),\n  @(\"21\u00d720=\", \"13\u00d763=\"),\n  @(\"76\u00d792=\", \"69\u00d725=\"),\n  @(\"98\u00d732=\", \"77\u00d725=\"),\n  @(\"88\u00d772=\", \"66\u00d755=\")\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Text = $old\n  $find.Replacement.ClearFormatting()\n  $find.Replacement.Text = $new\n  $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
